# Weekly update: a new price record for the week of 2022-01-28 (serial 44589)
# is inserted at row 50, pushing all existing records (previously rows 50-100)
# down by one row (new rows 51-101). Row 49 and above are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50 - this shifts rows 50..100 down to 51..101
# and the sheet dimension grows from A1:R100 to A1:R101 automatically.
$ws.Rows(50).Insert()

# Populate the new row 50 with this week's record.
$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 44589
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 100112022
$ws.Range("G50").Value = "Arveja Verde"
$ws.Range("H50").Value = "Perfection"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 52
$ws.Range("K50").Value = 23000
$ws.Range("L50").Value = 25000
$ws.Range("M50").Value = 24000
$ws.Range("N50").Value = "$/saco 25 kilos"
$ws.Range("O50").Value = "Carahue"
$ws.Range("P50").Value = 960
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
